# Applies the commit "add Tablet view and update WORK":
#  1. Moves the (hidden) "_GoBack" bookmark from the very last paragraph of
#     the document up to the empty paragraph right before the "BIO - HOME"
#     navigation line (this is simply where Word leaves the bookmark after
#     the most recent edit happened there).
#  2. Resizes the WORK page picture (Picture 15) taller (its width is kept
#     the same), reflecting the updated WORK/Tablet-view image.

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -------------------------------
$targetPara = $d.Paragraphs.Item(26)
$d.Bookmarks.Add("_GoBack", $targetPara.Range) | Out-Null

# --- 2. Resize the WORK picture ----------------------------------------
# (Width must stay fixed at 2073910 EMU / 163.3pt, so the aspect-ratio
#  lock is released for the resize and then restored to match the
#  original markup.)
$pic = $d.Shapes.Item(9)
$pic.LockAspectRatio = 0
$pic.Height = 427.5
$pic.LockAspectRatio = 1
